$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.556565037682828
$ws.Range("C2").Value = 0.612606004275778
$ws.Range("L2").Value = 0.599020579536134

$ws.Range("B3").Value = 0.503247309719211
$ws.Range("L3").Value = 0.605713809935328

$ws.Range("B4").Value = 0.680815013747804
$ws.Range("I4").Value = 0.68337742404715
$ws.Range("L4").Value = 0.689906588508866

$ws.Range("B5").Value = 0.704960018034767
$ws.Range("L5").Value = 0.777883926828007

$ws.Range("B6").Value = 0.410629863862209
$ws.Range("L6").Value = 0.389828757171604
